$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5813.06222222222
$ws.Range("D2").Value = 4865.03219436376
$ws.Range("E2").Value = 6761.09225008068
$ws.Range("F2").Value = 55.5459009356584
$ws.Range("G2").Value = 45.2685916647577
$ws.Range("H2").Value = 66.550298454881
$ws.Range("C3").Value = 7163.94964028777
$ws.Range("D3").Value = 3661.73651140874
$ws.Range("E3").Value = 10666.1627691668
$ws.Range("F3").Value = 48.8296604619294
$ws.Range("G3").Value = 23.4620529407042
$ws.Range("H3").Value = 79.4095214328846
$ws.Range("C4").Value = 6961.80842911877
$ws.Range("D4").Value = 4218.81038829578
$ws.Range("E4").Value = 9704.80646994177
$ws.Range("F4").Value = 66.9808285411527
$ws.Range("G4").Value = 40.341880927658
$ws.Range("H4").Value = 98.6762391667138
$ws.Range("C5").Value = 5812.71428571429
$ws.Range("D5").Value = -1414.1352472354
$ws.Range("E5").Value = 13039.563818664
$ws.Range("F5").Value = 53.6913057921562
$ws.Range("G5").Value = -1.72314443495992
$ws.Range("H5").Value = 140.351783136423
$ws.Range("C6").Value = 4704.07254901961
$ws.Range("D6").Value = 2954.8283999128
$ws.Range("E6").Value = 6453.31669812642
$ws.Range("F6").Value = 38.5442042766873
$ws.Range("G6").Value = 24.0547889416789
$ws.Range("H6").Value = 54.7259618303349
$ws.Range("C7").Value = 3318.10104529617
$ws.Range("D7").Value = 806.884317947733
$ws.Range("E7").Value = 5829.3177726446
$ws.Range("F7").Value = 27.7276273923899
$ws.Range("G7").Value = 10.8662674835068
$ws.Range("H7").Value = 47.1533873160854
$ws.Range("C8").Value = 6229.82608695652
$ws.Range("D8").Value = 3251.70509578649
$ws.Range("E8").Value = 9207.94707812656
$ws.Range("F8").Value = 47.8107071026915
$ws.Range("G8").Value = 23.7498354063453
$ws.Range("H8").Value = 76.5497712579372
$ws.Range("C9").Value = 5857.81606765328
$ws.Range("D9").Value = 4241.13842422649
$ws.Range("E9").Value = 7474.49371108006
$ws.Range("F9").Value = 54.1395423190276
$ws.Range("G9").Value = 38.5706200315785
$ws.Range("H9").Value = 71.4576906771792
$ws.Range("C10").Value = 7595.79189686924
$ws.Range("D10").Value = 5859.87431695783
$ws.Range("E10").Value = 9331.70947678066
$ws.Range("F10").Value = 61.6515614377799
$ws.Range("G10").Value = 45.1272559603083
$ws.Range("H10").Value = 80.0573375577302
$ws.Range("C11").Value = 10112.0432220039
$ws.Range("D11").Value = 8101.47413503325
$ws.Range("E11").Value = 12122.6123089746
$ws.Range("F11").Value = 91.9420953271509
$ws.Range("G11").Value = 70.3431066852557
$ws.Range("H11").Value = 116.27977013857
$ws.Range("C12").Value = 4501.90441176471
$ws.Range("D12").Value = 3255.31231883054
$ws.Range("E12").Value = 5748.49650469887
$ws.Range("F12").Value = 46.4901926570979
$ws.Range("G12").Value = 33.9369515661543
$ws.Range("H12").Value = 60.2199862979143
$ws.Range("C13").Value = 9529.77134146341
$ws.Range("D13").Value = 7149.09627142381
$ws.Range("E13").Value = 11910.446411503
$ws.Range("F13").Value = 91.4532013018841
$ws.Range("G13").Value = 65.0627731292905
$ws.Range("H13").Value = 122.062961828644
$ws.Range("C14").Value = 5722.04977375566
$ws.Range("D14").Value = 3045.41700719008
$ws.Range("E14").Value = 8398.68254032123
$ws.Range("F14").Value = 55.3556206788781
$ws.Range("G14").Value = 32.7393772633072
$ws.Range("H14").Value = 81.8252381028092
